$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Activate()

$ws.Range("A17").Value = "species"
$ws.Range("B17").Value = "Saccharomyces cerevisiae"
$ws.Range("A18").Value = "taxon_id"
$ws.Range("B18").Value = 559292

$ws.Range("A17:C18").Font.Size = 10
$ws.Range("A17:C18").Font.Name = "Arial"
$ws.Range("A17:C18").Font.Color = 0

$ws.Range("A17:C18").Select()
